$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Several match rows got re-ordered (their betting-odds payload rotated
#    down by one position within a small block of rows that share the same
#    match date / opening-odds timestamp). Columns A (index) and E (date)
#    stay put per row; everything from F..V (excluding the unchanged opening
#    timestamps K/O/S, which are identical across each block anyway) rotates.
#    We simply read the current F..V values for every row in a block, then
#    write them back shifted by one: new(row[i]) = old(row[i-1]), wrapping
#    the first row of the block around to the last row's data.
# ---------------------------------------------------------------------------

$rotCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$blocks = @( @(26,27,28), @(30,31,32), @(56,57,58), @(89,90) )

foreach ($rows in $blocks) {
    $data = @()
    foreach ($r in $rows) {
        $rowvals = @()
        foreach ($c in $rotCols) {
            $rowvals += , $ws.Range("$c$r").Value2
        }
        $data += , $rowvals
    }

    $n = $rows.Length
    for ($i = 0; $i -lt $n; $i++) {
        $srcIdx = ($i - 1 + $n) % $n
        $destRow = $rows[$i]
        $srcVals = $data[$srcIdx]
        for ($j = 0; $j -lt $rotCols.Length; $j++) {
            $ws.Range("$($rotCols[$j])$destRow").Value2 = $srcVals[$j]
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Four brand-new match rows were appended at the bottom (rows 110-113),
#    extending the used range from A1:V109 to A1:V113. Copy the formatting
#    of the last existing data row (109) down across the new rows first,
#    then fill in the values.
# ---------------------------------------------------------------------------

$ws.Range("A109:V109").Copy()
$ws.Range("A110:V113").PasteSpecial(-4122)

$newRows = @(
    @{A=109; E=45239.66666666666; F="Al Khaleej"; G=3; H="Al Taee";    I=1; J=1.96; K="04/11/2023 19:13"; L=2.12; M="09/11/2023 15:48"; N=3.46; O="04/11/2023 19:13"; P=3.75; Q="09/11/2023 15:56"; R=3.64; S="04/11/2023 19:13"; T=3.27; U="09/11/2023 15:48"; V="https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-khaleej-al-taee/WEHHIX55/"},
    @{A=110; E=45239.66666666666; F="Al Raed";    G=2; H="Al Shabab";  I=1; J=3.61; K="04/11/2023 19:13"; L=3.28; M="09/11/2023 15:53"; N=3.85; O="04/11/2023 19:13"; P=3.63; Q="09/11/2023 15:56"; R=1.95; S="04/11/2023 19:13"; T=2.15; U="09/11/2023 15:53"; V="https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-shabab/KxZgaq5p/"},
    @{A=111; E=45239.79166666666; F="Al Riyadh";  G=1; H="Al Fateh";   I=1; J=4.45; K="05/11/2023 19:12"; L=4.56; M="09/11/2023 18:58"; N=4.11; O="05/11/2023 19:12"; P=4.22; Q="09/11/2023 18:58"; R=1.72; S="05/11/2023 19:12"; T=1.69; U="09/11/2023 18:58"; V="https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-fateh/WCxbb3Kj/"},
    @{A=112; E=45239.79166666666; F="Damac";      G=2; H="Al Ahli SC"; I=2; J=4.18; K="06/11/2023 03:42"; L=4.06; M="09/11/2023 18:50"; N=4.22; O="06/11/2023 03:42"; P=3.95; Q="09/11/2023 18:50"; R=1.74; S="06/11/2023 03:42"; T=1.83; U="09/11/2023 18:50"; V="https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/damac-al-ahli-sc/vuLPGBzI/"}
)

$startRow = 110
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value2 = $row.A
    $ws.Range("B$r").Value2 = "saudi-arabia"
    $ws.Range("C$r").Value2 = "saudi-professional-league"
    $ws.Range("D$r").Value2 = "2023-2024"
    $ws.Range("E$r").Value2 = $row.E
    $ws.Range("F$r").Value2 = $row.F
    $ws.Range("G$r").Value2 = $row.G
    $ws.Range("H$r").Value2 = $row.H
    $ws.Range("I$r").Value2 = $row.I
    $ws.Range("J$r").Value2 = $row.J
    $ws.Range("K$r").Value2 = $row.K
    $ws.Range("L$r").Value2 = $row.L
    $ws.Range("M$r").Value2 = $row.M
    $ws.Range("N$r").Value2 = $row.N
    $ws.Range("O$r").Value2 = $row.O
    $ws.Range("P$r").Value2 = $row.P
    $ws.Range("Q$r").Value2 = $row.Q
    $ws.Range("R$r").Value2 = $row.R
    $ws.Range("S$r").Value2 = $row.S
    $ws.Range("T$r").Value2 = $row.T
    $ws.Range("U$r").Value2 = $row.U
    $ws.Range("V$r").Value2 = $row.V
}
